# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get new values that look like plain decimal numbers
# (e.g. "302.91"); mark them as Text first so Excel keeps them as the
# original string type instead of silently coercing them to numbers.
# (Multi-area "A1,A2,..." range strings only apply to the first area, so
# each cell is set individually.)
$textCells = "D5","D6","D9","D10","D14","D17","D19","D20","D22","D23","D24","D27","D30","D32","D33","D34","D35","D36","D38","D41","D44","D49","D50","D51"
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.890.13"
$ws.Range("E2").Value = "  -1.35%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.338.74"
$ws.Range("E3").Value = "  +0.22%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "302.91"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6 - Solana
$ws.Range("D6").Value = "93.99"
$ws.Range("E6").Value = "  -4.50%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.44%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  -1.98%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "34.03"
$ws.Range("E10").Value = "  -4.89%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -2.25%  "

# Row 12 - Chainlink
$ws.Range("E12").Value = "  -3.88%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.10%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.73"
$ws.Range("E14").Value = "  -2.93%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.702.05"
$ws.Range("E15").Value = "  +0.26%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.357.54"
$ws.Range("E16").Value = "  +1.82%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.794"
$ws.Range("E17").Value = "  +0.01%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.834.64"
$ws.Range("E18").Value = "  -1.23%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").Value = "12.03"
$ws.Range("E19").Value = "  -6.04%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").Value = "  +1.61%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  -1.71%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "67.84"
$ws.Range("E22").Value = "  -0.48%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "235.25"
$ws.Range("E23").Value = "  -1.09%  "

# Row 24 - ImmutableX
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  -2.41%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.27%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  -1.30%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "24.60"
$ws.Range("E27").Value = "  -2.14%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  -0.61%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  -0.42%  "

# Row 30 - InjectiveProtocol
$ws.Range("D30").Value = "31.31"
$ws.Range("E30").Value = "  -6.29%  "

# Row 31 - FirstDigitalUSD
$ws.Range("E31").Value = "  +0.03%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "4.98"
$ws.Range("E32").Value = "  -0.85%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.0743"
$ws.Range("E33").Value = "  +5.05%  "

# Row 34 - Celestia
$ws.Range("D34").Value = "17.25"
$ws.Range("E34").Value = "  -3.81%  "

# Row 35 - RenderToken
$ws.Range("D35").Value = "4.39"
$ws.Range("E35").Value = "  -2.20%  "

# Row 36 - WEMIXToken
$ws.Range("D36").Value = "2.33"
$ws.Range("E36").Value = "  -0.88%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  +1.56%  "

# Row 38 - Monero
$ws.Range("D38").Value = "125.32"
$ws.Range("E38").Value = "  -23.72%  "

# Row 40 - LidoDAOToken
$ws.Range("E40").Value = "  -1.12%  "

# Row 41 - EnergySwap
$ws.Range("D41").Value = "22.35"
$ws.Range("E41").Value = "  +17.06%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  -1.56%  "

# Row 43 - Maker
$ws.Range("D43").Value = "1.936.85"
$ws.Range("E43").Value = "  -2.73%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "0.0281"
$ws.Range("E44").Value = "  -0.36%  "

# Row 45 - FraxShare
$ws.Range("E45").Value = "  -4.89%  "

# Row 46 - ApeXProtocol
$ws.Range("E46").Value = "  +0.85%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  -3.66%  "

# Row 48/49 - HuobiToken and RocketPoolETH swap rank order
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.567.95"
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "2.87"
$ws.Range("E49").Value = "  +0.07%  "

# Row 50 - MultiversX
$ws.Range("D50").Value = "52.71"
$ws.Range("E50").Value = "  -2.89%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "71.49"
$ws.Range("E51").Value = "  -2.16%  "
